$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting the old
# "Late"/"heading"/"Outstanding" columns (N,O,P) one to the right (O,P,Q).
$ws.Columns("N").Insert()

# The new column inherits no explicit width from the insert; give it the
# same width as column M (10.7109375 in OOXML units - the nearest
# representable width via ColumnWidth is used).
$ws.Columns("N").ColumnWidth = 9.83

# The Insert() shift re-serialises the moved "Outstanding" values (old
# column P, now Q) through floating point, introducing tiny binary noise
# (e.g. 963.77 -> 963.76999999999998). Re-apply the exact literal values
# so the saved numbers match the originals precisely.
$ws.Range("Q3").Value = 0
$ws.Range("Q4").Value = 963.77
$ws.Range("Q5").Value = 963.77
$ws.Range("Q6").Value = 963.77
$ws.Range("Q7").Value = 963.77
$ws.Range("Q8").Value = 963.77
$ws.Range("Q9").Value = 963.77
$ws.Range("Q10").Value = 963.77
$ws.Range("Q11").Value = 963.77
$ws.Range("Q12").Value = 963.77
$ws.Range("Q13").Value = 1079.81

# Make "Repayment schedule" the active/selected sheet with S10 selected.
$ws.Activate() | Out-Null
$ws.Range("S10").Select() | Out-Null
